# Generate Report for Handoff
# - Updates the "Latest Handoff Datetime" (zh-cn / de-de sheets) and the
#   corresponding "Latest HO Xliff Generate Date" (Overview sheet) for the
#   file 4c2c0fec-d838-4153-9d5b-bfcf2a00f122.md, which now has a handoff
#   priority ("ht") recorded for rows 8-13 on both language sheets.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" column (G), rows 8-13
$wsOverview.Range("G8:G13").Value = "2016-08-25 20:19:51"

# zh-cn sheet: "Latest Handoff Datetime" column (H), rows 8-13
$wsZhCn.Range("H8:H13").Value = "2016-08-25 20:19:46"

# de-de sheet: "Latest Handoff Datetime" column (H), rows 8-13
$wsDeDe.Range("H8:H13").Value = "2016-08-25 20:19:51"

# zh-cn sheet: "Priority" column (E), rows 8-13 now show handoff type "ht"
$wsZhCn.Range("E8:E13").Value = "ht"

# de-de sheet: "Priority" column (E), rows 8-13 now show handoff type "ht"
$wsDeDe.Range("E8:E13").Value = "ht"
